# Weekly data update: append the new week's row (2018-06-14, serial 43265)
# to both "Sheet1" and "underReview", matching the source workbook's
# pattern of carrying the date's number format forward from the prior row.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("underReview")

# ---- Sheet1: new row 8 ----
# Copy the date-formatted style from the previous row's date cell so the
# new date cell reuses the existing "short date" style instead of minting
# a new one.
$ws1.Range("A7").Copy()
$ws1.Range("A8").PasteSpecial(-4122)   # xlPasteFormats

$ws1.Range("A8").Value = 43265
$ws1.Range("B8").Value = 1145
$ws1.Range("C8").Value = 128
$ws1.Range("D8").Value = 72
$ws1.Range("E8").Value = 945
$ws1.Range("F8").Value = 2
$ws1.Range("G8").Value = 34
$ws1.Range("H8").Value = 26
$ws1.Range("I8").Value = 7
$ws1.Range("J8").Value = 17
$ws1.Range("K8").Value = 3
$ws1.Range("L8").Value = 36
$ws1.Range("M8").Value = 3

# ---- underReview: new row 3 ----
$ws2.Range("A2").Copy()
$ws2.Range("A3").PasteSpecial(-4122)   # xlPasteFormats

$ws2.Range("A3").Value = 43265
$ws2.Range("B3").Value = 53
$ws2.Range("C3").Value = 0
$ws2.Range("D3").Value = 18
$ws2.Range("E3").Value = "NA"
$ws2.Range("F3").Value = 2
$ws2.Range("G3").Value = 12
$ws2.Range("H3").Value = 3
$ws2.Range("I3").Value = 17
$ws2.Range("J3").Value = 1

# Column A on underReview now holds dates like Sheet1's column A -- widen
# it to fit, mirroring Sheet1's existing column setup.
$ws2.Columns.Item(1).AutoFit()

# ---- Selections / active tab ----
# underReview is the active sheet in the source file; move its selection
# down to the newly added row...
$ws2.Range("A3").Select() | Out-Null

# ...then switch focus to Sheet1 (making it the active/saved tab) and
# leave the selection on B12, matching the saved view state.
$ws1.Activate() | Out-Null
$ws1.Range("B12").Select() | Out-Null
